$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update recalculated values in column F (total column) for several rows.
$ws.Range("F5").Value = 378091.86347555893
$ws.Range("F15").Value = 354171.50838907313
$ws.Range("F20").Value = 350683.44765434955
$ws.Range("F25").Value = 361684.42839035933
$ws.Range("F30").Value = 344586.54493832408
$ws.Range("F32").Value = 1918.3988340409351
$ws.Range("F33").Value = 16942.006527913905
$ws.Range("F34").Value = 22442.227475102762
$ws.Range("F35").Value = 371858.81673632929
$ws.Range("F36").Value = 120572.562235836
$ws.Range("F40").Value = 376095.50883206091
$ws.Range("F45").Value = 366767.9169788835
$ws.Range("F50").Value = 385884.80354058166

# F6's number-format style got cleared back to the default/"Normal" style.
$ws.Range("F6").Style = "Normal"
$ws.Range("F6").Value = 123661.77577379883

# Move the selection to D7 (also resets the scrolled topLeftCell back to default).
[void]$ws.Range("D7").Select()
